$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 15) lists MIGS.ba attribute names starting at column A.
# We are inserting a new attribute "culture_collection" right before
# "encoded_traits", which currently lives in column R (18th column).
# That pushes "encoded_traits" and everything after it one column to the
# right (R -> S, S -> T, ... AG -> AH).
#
# Excel's Columns.Insert() shifts the cell VALUES/styles correctly, but the
# cell Comments in this runtime stay anchored to their original column, so
# we capture the old comment text first and re-apply it one column to the
# right ourselves.

$insertCol = 18  # column R

# 1) Capture existing comments for columns R (18) .. AG (33) before the shift.
$oldComments = @{}
for ($i = $insertCol; $i -le 33; $i++) {
    $cell = $ws.Cells.Item(15, $i)
    $cm = $cell.Comment
    if ($cm) {
        $oldComments[$i] = $cm.Text()
    }
}

# 2) Insert a new column at R, shifting encoded_traits..trophic_level right.
$ws.Columns.Item($insertCol).Insert()

# 3) Re-home the captured comments one column to the right, starting from
#    the far end so we never overwrite a cell before reading it.
for ($i = 33; $i -ge $insertCol; $i--) {
    $destCol = $i + 1
    $destCell = $ws.Cells.Item(15, $destCol)
    $destCell.ClearComments()
    if ($oldComments.ContainsKey($i)) {
        $destCell.AddComment($oldComments[$i])
    }
}

# 4) Populate the new culture_collection column header + comment.
$newCell = $ws.Cells.Item(15, $insertCol)
$newCell.Value() = "culture_collection"
$newCell.ClearComments()
$newCell.AddComment("Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier")
